$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "delivered"
$ws.Range("D2").Value = "SEVILLA - SPAIN"
$ws.Range("E2").Value = "MUNICH - GERMANY"
$ws.Range("F2").Value = "{'timestamp': '2024-02-08T08:28:00', 'location': {'address': {'addressLocality': 'MUNICH - GERMANY'}}, 'description': 'Delivered', 'pieceIds': ['JD014600011333431260']}"

$ws.Range("F1").Value = "Events (newest)"

$ws.Columns.Item(4).ColumnWidth = 23.5703125
$ws.Columns.Item(5).ColumnWidth = 30.7109375
$ws.Columns.Item(6).ColumnWidth = 16.85546875

$ws.Range("G9").Select()
